$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column B. This shifts the existing
# B:F columns ("a","b","c","d","message" and their data) two columns to
# the right (now D:H), matching the diff which shows the header/data for
# those columns unchanged in content but moved to columns D-H.
$ws.Range("B1:C1").EntireColumn.Insert()

# Give the two new header cells (B1/C1) the same style as the rest of the
# header row (bordered/bold/centered style used by D1:H1), then set their
# text.
$ws.Range("D1").Copy()
$ws.Range("B1:C1").PasteSpecial(-4122)
$ws.Range("B1").Value = "Unnamed: 0"
$ws.Range("C1").Value = "Unnamed: 0.1"

# The new B/C data columns mirror column A (the old DataFrame index),
# duplicated from a double reset_index(). Column A keeps its special
# (bordered) style, but the new B/C data cells should have no explicit
# style, matching the other plain data cells.
$ws.Range("B2:C4").ClearFormats()

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 2
